# Update cryptos list with latest price/volume figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.857.37'
$ws.Range("E2").Value = '  +1.55%  '
$ws.Range("D3").Value = '1.769.02'
$ws.Range("E3").Value = '  +1.94%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '327.92'
$ws.Range("E5").Value = '  +2.11%  '
$ws.Range("E7").Value = '  -2.81%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3561'
$ws.Range("E8").Value = '  +1.23%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07445'
$ws.Range("E9").Value = '  +1.38%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.02'
$ws.Range("E10").Value = '  +1.23%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.098'
$ws.Range("E11").Value = '  +2.06%  '
$ws.Range("E12").Value = '  -0.28%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.96'
$ws.Range("E13").Value = '  +3.09%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.025'
$ws.Range("E14").Value = '  +2.06%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.250'
$ws.Range("E15").Value = '  +2.98%  '
$ws.Range("D16").Value = '1.773.01'
$ws.Range("E16").Value = '  +1.86%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '93.24'
$ws.Range("E17").Value = '  +2.68%  '
$ws.Range("E18").Value = '  +0.97%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06438'
$ws.Range("E19").Value = '  +1.49%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  -0.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.13'
$ws.Range("E21").Value = '  +3.20%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.780'
$ws.Range("E22").Value = '  +1.07%  '
$ws.Range("D23").Value = '27.914.26'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.29'
$ws.Range("E24").Value = '  +2.20%  '
$ws.Range("E25").Value = '  +0.24%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.07'
$ws.Range("E26").Value = '  +0.62%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.35'
$ws.Range("E27").Value = '  +2.93%  '
$ws.Range("D28").Value = '1.981.55'
$ws.Range("E28").Value = '  +2.25%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.152'
$ws.Range("E29").Value = '  +5.90%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.90'
$ws.Range("E31").Value = '  +6.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09183'
$ws.Range("E32").Value = '  +0.76%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.616'
$ws.Range("E33").Value = '  +4.40%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.653'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.87'
$ws.Range("E35").Value = '  +2.78%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02292'
$ws.Range("E36").Value = '  +1.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06087'
$ws.Range("E37").Value = '  +1.84%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2102'
$ws.Range("E38").Value = '  +2.14%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6328'
$ws.Range("E39").Value = '  +1.80%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.959'
$ws.Range("E40").Value = '  +1.77%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.184'
$ws.Range("E41").Value = '  +0.72%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.394'
$ws.Range("E42").Value = '  +1.62%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.891'
$ws.Range("E43").Value = '  +2.47%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.38'
$ws.Range("E44").Value = '  +2.60%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.742'
$ws.Range("E45").Value = '  +1.37%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5910'
$ws.Range("E46").Value = '  +2.19%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '122.33'
$ws.Range("E47").Value = '  +0.54%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.958'
$ws.Range("E48").Value = '  +2.37%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06904'
$ws.Range("E49").Value = '  +1.15%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.137'
$ws.Range("E50").Value = '  +2.61%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.01'
$ws.Range("E51").Value = '  +2.77%  '
